$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D: "Done" / "-" status for rows 1-4 (existing OS PPT entries in column C)
$ws.Range("D1").Value2 = "Done"
$ws.Range("D2").Value2 = "Done"
$ws.Range("D3").Value2 = "-"
$ws.Range("D4").Value2 = "Done"

# Column C: continue the OS / PDD PPT list for rows 5-6
$ws.Range("C5").Value2 = "OS 6"
$ws.Range("C6").Value2 = "PDD 4"

# Note about recent topics not covered
$ws.Range("E2").Value2 = "Recent topics not done"

# Column C: OODP topic ranges for rows 13-16
$ws.Range("C13").Value2 = "1, 2"
$ws.Range("C14").Value2 = "3, 4"
$ws.Range("C15").Value2 = "5, 6"
$ws.Range("C16").Value2 = "All"

# Update the active selection to C7
[void]$ws.Range("C7").Select()
